$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.162.19'
$ws.Range("E2").Value = '  -4.39%  '
$ws.Range("D3").Value = '1.656.28'
$ws.Range("E3").Value = '  -3.06%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5175'
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06432'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2569'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.89'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07776'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.06%  '
$ws.Range("D12").Value = '1.662.23'
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("D13").Value = '1.884.49'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("E14").Value = '  -5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5536'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.08%  '
$ws.Range("D16").Value = '0.0₅8053'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.38'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.95%  '
$ws.Range("D18").Value = '26.191.44'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.43'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.383'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.908'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.10'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.761'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1162'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.973'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.76'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05275'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.367'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.228'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.574'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.762'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.363'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9242'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.55%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5722'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.167.67'
$ws.Range("E39").Value = '  +11.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01593'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.69%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8369'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.666'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.95'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = '1.794.96'
$ws.Range("E45").Value = '  -3.02%  '
$ws.Range("E46").Value = '  -6.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4505'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.93'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.008'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.911'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05080'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.92%  '
